# Add a new entry (row 33) to the working-hours log on Sheet1, mirroring
# the formatting of the most recent existing entry (row 32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy cell formatting from the row above (row 32) for the columns that
# carry an explicit style (date column A, day-of-week column B, wrapped
# notes column E). This keeps number formats / wrap formatting identical
# to the rest of the log without dragging along unrelated columns.
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B32").Copy()
$ws.Range("B33").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("E32").Copy()
$ws.Range("E33").PasteSpecial(-4122)   # xlPasteFormats

# New log entry: Tuesday, June 11 2024, 8 hours.
$ws.Range("A33").Value = 45454
$ws.Range("B33").Value = "T"
$ws.Range("C33").Value = 8

# Write the TODO/"Links" column note first and the Notes column second so
# the shared-string table picks up the two new strings in the same order
# as the source workbook.
$ws.Range("G33").Value = "want to run grid search on xgboost"
$ws.Range("E33").Value = "Meeting, using quantile regression for error bars. Xgboost is working, grid search on xgboost"

# Match the row height used by the other wrapped two-line entries.
$ws.Rows.Item(33).RowHeight = 28.5
